# Mise a jour du fichier Excel via Shiny
# Appends new forecast/calibration quarters (2022-01-01 .. 2025-10-01) and
# the corresponding yearly calibration rows (2022 .. 2025) across the
# INDI, ETALONNAGE, PREVISION and VATRIM sheets.

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-NumberCell($ws, $row, $col, $num) {
    $ws.Cells.Item($row, $col).Value = $num
}

# ---------------------------------------------------------------------
# INDI sheet: append quarterly rows 106-121 (columns A=DATE, B=INDI_NB, C=IND_base)
# ---------------------------------------------------------------------
$wsIndi = $wb.Worksheets.Item("INDI")

$indiData = @(
    @("2022-01-01", 1086.56222288586, 148.911907134174),
    @("2022-04-01", 1128.67743752971, 154.683741272972),
    @("2022-07-01", 1121.66342401112, 153.72248005137),
    @("2022-10-01", 1163.61090038726, 159.47132588373),
    @("2023-01-01", 1134.62400164452, 155.498709973872),
    @("2023-04-01", 1173.42592840471, 160.81646241606),
    @("2023-07-01", 1167.84688534618, 160.051861986991),
    @("2023-10-01", 1206.49728546378, 165.348847904396),
    @("2024-01-01", 1181.06900006139, 161.863935218618),
    @("2024-04-01", 1216.87786271158, 166.771492206364),
    @("2024-07-01", 1212.59510260031, 166.184545630702),
    @("2024-10-01", 1248.26708340636, 171.073343143811),
    @("2025-01-01", 1226.05347727255, 168.028997975131),
    @("2025-04-01", 1259.15855049025, 172.566004218151),
    @("2025-07-01", 1256.04678973131, 172.139541545883),
    @("2025-10-01", 1289.02821050788, 176.659601386348)
)

$row = 106
foreach ($r in $indiData) {
    Set-TextCell $wsIndi $row 1 $r[0]
    Set-NumberCell $wsIndi $row 2 $r[1]
    Set-NumberCell $wsIndi $row 3 $r[2]
    $row++
}

# ---------------------------------------------------------------------
# ETALONNAGE sheet: append yearly rows 28-31
# (columns A=ANNEE, B=PROD_CST, C=Indicateur_Base, D=evol_Prod,
#  E=evol_indicateur, F=decision)
# ---------------------------------------------------------------------
$wsEtal = $wb.Worksheets.Item("ETALONNAGE")

$etalData = @(
    @("2022", 2838773.95539881, 154.197363585561, 4.39850350114608, 4.35771211321176, "Acceptable"),
    @("2023", 2954526.51451073, 160.42897057033, 4.07755463909969, 4.04131876178953, "Acceptable"),
    @("2024", 3066800.91883488, 166.473329049874, 3.80008112205918, 3.76762280407117, "Acceptable"),
    @("2025", 3175933.32625752, 172.348536281378, 3.55850967542084, 3.52921832286057, "Acceptable")
)

$row = 28
foreach ($r in $etalData) {
    Set-TextCell $wsEtal $row 1 $r[0]
    Set-NumberCell $wsEtal $row 2 $r[1]
    Set-NumberCell $wsEtal $row 3 $r[2]
    Set-NumberCell $wsEtal $row 4 $r[3]
    Set-NumberCell $wsEtal $row 5 $r[4]
    Set-TextCell $wsEtal $row 6 $r[5]
    $row++
}

# ---------------------------------------------------------------------
# PREVISION sheet: append quarterly rows 102-117
# (columns A=time, B=Prod_brut, C=VA_brut_M1, D=VA_brut_M2, E=CI_brut_M2,
#  F=IND_trim, G=Ut_1, H=RI, I=Var_M1, J=Var_M2, K=verif)
# ---------------------------------------------------------------------
$wsPrev = $wb.Worksheets.Item("PREVISION")

$prevData = @(
    @("2022-01-01", 685149.053047744, 692555.335208627, 685149.053047744, 0, 37.2279767835435, -146501.204705833, 18603.0881891698, 4.80862566808695, 4.85545451130682, "Acceptable"),
    @("2022-04-01", 711952.114199086, 686889.716760687, 711952.114199086, 0, 38.6709353182431, -146501.204705833, 17762.4283226645, 4.26642347843802, 4.30619844178988, "Acceptable"),
    @("2022-07-01", 707488.239686158, 669690.251591192, 707488.239686158, 0, 38.4306200128425, -146501.204705833, 17425.9549121872, 4.44851971939908, 4.49032767009454, "Acceptable"),
    @("2022-10-01", 734184.548465822, 689654.457620217, 734184.548465822, 0, 39.8678314709324, -146501.204705833, 17298.5194372321, 3.94130316255941, 3.97682190940287, "Acceptable"),
    @("2023-01-01", 715736.640122968, 723189.05373652, 715736.640122968, 0, 38.8746774934679, -146501.204705833, 18603.0881891698, 4.42328821547013, 4.46436975124769, "Acceptable"),
    @("2023-04-01", 740431.051027391, 714122.721692433, 740431.051027391, 0, 40.204115604015, -146501.204705833, 17762.4283226645, 3.96468374285379, 4.00011970753722, "Acceptable"),
    @("2023-07-01", 736880.423694059, 697264.132649229, 736880.423694059, 0, 40.0129654967477, -146501.204705833, 17425.9549121872, 4.11740815884962, 4.15444135452185, "Acceptable"),
    @("2023-10-01", 761478.399666311, 715072.56484953, 761478.399666311, 0, 41.3372119760989, -146501.204705833, 17298.5194372321, 3.68562936822341, 3.71757363424803, "Acceptable"),
    @("2024-01-01", 745295.272158177, 752792.265379507, 745295.272158177, 0, 40.4659838046546, -146501.204705833, 18603.0881891698, 4.09342639936763, 4.12981959818781, "Acceptable"),
    @("2024-04-01", 768084.830361531, 740566.669144836, 768084.830361531, 0, 41.6928730515911, -146501.204705833, 17762.4283226645, 3.70299762899755, 3.734821668509, "Acceptable"),
    @("2024-07-01", 765359.186383944, 723981.099815733, 765359.186383944, 0, 41.5461364076755, -146501.204705833, 17425.9549121872, 3.83168528474205, 3.86477395438438, "Acceptable"),
    @("2024-10-01", 788061.629931232, 739828.887891372, 788061.629931232, 0, 42.7683357859527, -146501.204705833, 17298.5194372321, 3.46207143984774, 3.49100253882053, "Acceptable"),
    @("2025-01-01", 773924.395949848, 781464.566917298, 773924.395949848, 0, 42.0072494937826, -146501.204705833, 18603.0881891698, 3.80879332272848, 3.84131294819148, "Acceptable"),
    @("2025-04-01", 794993.202303494, 766297.820213383, 794993.202303494, 0, 43.1415010545378, -146501.204705833, 17762.4283226645, 3.474521895275, 3.50330730126478, "Acceptable"),
    @("2025-07-01", 793012.808409853, 749923.972395782, 793012.808409853, 0, 43.0348853864706, -146501.204705833, 17425.9549121872, 3.58336323788724, 3.61315608643342, "Acceptable"),
    @("2025-10-01", 814002.919594324, 763987.387088852, 814002.919594324, 0, 44.1649003465869, -146501.204705833, 17298.5194372321, 3.26541712453752, 3.29178438307627, "Acceptable")
)

$row = 102
foreach ($r in $prevData) {
    Set-TextCell $wsPrev $row 1 $r[0]
    Set-NumberCell $wsPrev $row 2 $r[1]
    Set-NumberCell $wsPrev $row 3 $r[2]
    Set-NumberCell $wsPrev $row 4 $r[3]
    Set-NumberCell $wsPrev $row 5 $r[4]
    Set-NumberCell $wsPrev $row 6 $r[5]
    Set-NumberCell $wsPrev $row 7 $r[6]
    Set-NumberCell $wsPrev $row 8 $r[7]
    Set-NumberCell $wsPrev $row 9 $r[8]
    Set-NumberCell $wsPrev $row 10 $r[9]
    Set-TextCell $wsPrev $row 11 $r[10]
    $row++
}

# ---------------------------------------------------------------------
# VATRIM sheet: append quarterly rows 102-117 (columns A=time, B=VA_trim)
# ---------------------------------------------------------------------
$wsVatrim = $wb.Worksheets.Item("VATRIM")

$vatrimData = @(
    @("2022-01-01", 692555.335208627),
    @("2022-04-01", 686889.716760687),
    @("2022-07-01", 669690.251591192),
    @("2022-10-01", 689654.457620217),
    @("2023-01-01", 723189.05373652),
    @("2023-04-01", 714122.721692433),
    @("2023-07-01", 697264.132649229),
    @("2023-10-01", 715072.56484953),
    @("2024-01-01", 752792.265379507),
    @("2024-04-01", 740566.669144836),
    @("2024-07-01", 723981.099815733),
    @("2024-10-01", 739828.887891372),
    @("2025-01-01", 781464.566917298),
    @("2025-04-01", 766297.820213383),
    @("2025-07-01", 749923.972395782),
    @("2025-10-01", 763987.387088852)
)

$row = 102
foreach ($r in $vatrimData) {
    Set-TextCell $wsVatrim $row 1 $r[0]
    Set-NumberCell $wsVatrim $row 2 $r[1]
    $row++
}
